$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 149997.5
$ws.Range("J3").Value = 149997.5
$ws.Range("L3").Value = 149997.5
$ws.Range("N3").Value = -150225.5
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 14
$ws.Range("N4").Value = -728
$ws.Range("H17").Value = 687892.1
$ws.Range("J17").Value = 687892.1
$ws.Range("L17").Value = 2063676.3
$ws.Range("N17").Value = -2064012.3
$ws.Range("H18").Value = 1453.2222
$ws.Range("I18").Value = 1453.2222
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1453.2222
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1169.2222
$ws.Range("N18").ClearContents()
$ws.Range("H64").Value = 11001.462
$ws.Range("I64").Value = 4266.5
$ws.Range("J64").Value = 12226
$ws.Range("K64").Value = 4266.5
$ws.Range("L64").Value = 12226
$ws.Range("M64").Value = -4018.5
$ws.Range("N64").Value = -12722
$ws.Range("H67").Value = 11001.462
$ws.Range("I67").Value = 4266.5
$ws.Range("J67").Value = 12226
$ws.Range("K67").Value = 4266.5
$ws.Range("L67").Value = 12226
$ws.Range("M67").Value = -3408.5
$ws.Range("N67").Value = -13942
$ws.Range("H80").Value = 1289.75
$ws.Range("I80").Value = 1309.25
$ws.Range("J80").Value = 1280
$ws.Range("K80").Value = 3927.75
$ws.Range("L80").Value = 3840
$ws.Range("M80").Value = -2929.75
$ws.Range("N80").Value = -5836
$ws.Range("H83").Value = 1289.75
$ws.Range("I83").Value = 1309.25
$ws.Range("J83").Value = 1280
$ws.Range("K83").Value = 11783.25
$ws.Range("L83").Value = 11520
$ws.Range("M83").Value = -6791.25
$ws.Range("N83").Value = -21504
$ws.Range("H87").Value = 70833.336
$ws.Range("J87").Value = 71000
$ws.Range("L87").Value = 71000
$ws.Range("N87").Value = -73496
$ws.Range("H90").Value = 70833.336
$ws.Range("J90").Value = 71000
$ws.Range("L90").Value = 213000
$ws.Range("N90").Value = -225480
$ws.Range("H97").Value = 166666
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H100").Value = 1092.8108
$ws.Range("I100").Value = 1002.4483
$ws.Range("J100").Value = 1420.375
$ws.Range("K100").Value = 1002.4483
$ws.Range("L100").Value = 1420.375
$ws.Range("M100").Value = -461.4483
$ws.Range("N100").Value = -2502.375
$ws.Range("H102").Value = 149997.5
$ws.Range("J102").Value = 149997.5
$ws.Range("L102").Value = 149997.5
$ws.Range("N102").Value = -156487.5
$ws.Range("H105").Value = 48666.668
$ws.Range("J105").Value = 48666.668
$ws.Range("L105").Value = 48666.668
$ws.Range("N105").Value = -55654.668
$ws.Range("H132").Value = 1943.1025
$ws.Range("I132").Value = 1580.7097
$ws.Range("K132").Value = 4742.1291
$ws.Range("M132").Value = -2212.1291
$ws.Range("H141").Value = 2194.9285
$ws.Range("I141").Value = 1979.1538
$ws.Range("K141").Value = 5937.4614
$ws.Range("M141").Value = -757.4614000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 97.59999999999999
$ws.Range("I3").Value = 97.59999999999999
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 97.59999999999999
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 17.40000000000001
$ws.Range("N3").ClearContents()
$ws.Range("H132").Value = 2588.7036
$ws.Range("I132").Value = 2468.475
$ws.Range("K132").Value = 7405.424999999999
$ws.Range("M132").Value = -4875.424999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1836.8948
$ws.Range("I107").Value = 1185.1818
$ws.Range("K107").Value = 1185.1818
$ws.Range("M107").Value = 734.8181999999999
$ws.Range("H134").Value = 1985.425
$ws.Range("I134").Value = 2014.4286
$ws.Range("K134").Value = 6043.2858
$ws.Range("M134").Value = -3508.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 718791.1
$ws.Range("J3").Value = 9732.333000000001
$ws.Range("L3").Value = 9732.333000000001
$ws.Range("N3").Value = -9958.333000000001
$ws.Range("H7").Value = 479.6
$ws.Range("I7").Value = 479.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 479.6
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -366.6
$ws.Range("N7").ClearContents()
$ws.Range("H31").Value = 25889.441
$ws.Range("J31").Value = 3875.0667
$ws.Range("L31").Value = 3875.0667
$ws.Range("N31").Value = -4465.066699999999
$ws.Range("H34").Value = 25889.441
$ws.Range("J34").Value = 3875.0667
$ws.Range("L34").Value = 3875.0667
$ws.Range("N34").Value = -4279.066699999999
$ws.Range("H99").Value = 6024.4443
$ws.Range("I99").Value = 2844
$ws.Range("K99").Value = 2844
$ws.Range("M99").Value = -1346
$ws.Range("H105").Value = 1547.2
$ws.Range("I105").Value = 1434
$ws.Range("K105").Value = 1434
$ws.Range("M105").Value = 313
$ws.Range("H107").Value = 970.25
$ws.Range("I107").Value = 480.5
$ws.Range("J107").Value = 1460
$ws.Range("K107").Value = 480.5
$ws.Range("L107").Value = 1460
$ws.Range("M107").Value = 1439.5
$ws.Range("N107").Value = -5300
$ws.Range("H126").Value = 6024.4443
$ws.Range("I126").Value = 2844
$ws.Range("K126").Value = 8532
$ws.Range("M126").Value = -6062
$ws.Range("H132").Value = 4147.561
$ws.Range("I132").Value = 3989.7646
$ws.Range("J132").Value = 4914
$ws.Range("K132").Value = 11969.2938
$ws.Range("L132").Value = 14742
$ws.Range("M132").Value = -9439.293799999999
$ws.Range("N132").Value = -19802

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 669.6
$ws.Range("I7").Value = 574.5
$ws.Range("J7").Value = 733
$ws.Range("K7").Value = 1723.5
$ws.Range("L7").Value = 2199
$ws.Range("M7").Value = -1611.5
$ws.Range("N7").Value = -2423
$ws.Range("H34").Value = 2315931.5
$ws.Range("I34").Value = 67.666664
$ws.Range("K34").Value = 202.999992
$ws.Range("M34").Value = -118.999992
$ws.Range("H39").Value = 1698.9
$ws.Range("J39").Value = 1698.9
$ws.Range("L39").Value = 5096.700000000001
$ws.Range("N39").Value = -5684.700000000001
$ws.Range("H55").Value = 1164.375
$ws.Range("J55").Value = 1978.3334
$ws.Range("L55").Value = 5935.0002
$ws.Range("N55").Value = -6289.0002
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 41668524
$ws.Range("I102").Value = 1327.4445
$ws.Range("K102").Value = 1327.4445
$ws.Range("M102").Value = 294.5554999999999
$ws.Range("H126").Value = 17407.5
$ws.Range("I126").Value = 30378
$ws.Range("J126").Value = 4437
$ws.Range("K126").Value = 91134
$ws.Range("L126").Value = 13311
$ws.Range("M126").Value = -88664
$ws.Range("N126").Value = -18251
$ws.Range("H132").Value = 3310.0833
$ws.Range("I132").Value = 3143.1
$ws.Range("K132").Value = 9429.299999999999
$ws.Range("M132").Value = -6899.299999999999
$ws.Range("H135").Value = 49749.5
$ws.Range("J135").Value = 49749.5
$ws.Range("L135").Value = 49749.5
$ws.Range("N135").Value = -59889.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 3850
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 3850
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 3850
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -4440
$ws.Range("H46").Value = 1566.6666
$ws.Range("I46").Value = 1200
$ws.Range("K46").Value = 1200
$ws.Range("M46").Value = -1012
$ws.Range("H93").Value = 6000
$ws.Range("I93").Value = 6000
$ws.Range("K93").Value = 6000
$ws.Range("M93").Value = -4752
$ws.Range("H136").Value = 5295.6
$ws.Range("I136").Value = 4750.4375
$ws.Range("K136").Value = 14251.3125
$ws.Range("M136").Value = -11701.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16765.666
$ws.Range("J45").Value = 17405
$ws.Range("L45").Value = 17405
$ws.Range("N45").Value = -18387
$ws.Range("H70").Value = 30000
$ws.Range("I70").Value = 30000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685
$ws.Range("H73").Value = 30000
$ws.Range("I73").Value = 30000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908
$ws.Range("H80").Value = 34995
$ws.Range("J80").Value = 34995
$ws.Range("L80").Value = 34995
$ws.Range("N80").Value = -36991
$ws.Range("H83").Value = 34995
$ws.Range("J83").Value = 34995
$ws.Range("L83").Value = 104985
$ws.Range("N83").Value = -114969
